$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row before row 135, pushing existing rows 135-259 down to 136-260.
$ws.Rows(135).Insert()

# Populate the newly inserted row 135 with its data.
$ws.Range("A135").Value = 10
$ws.Range("B135").Value = "Vega Modelo de Temuco"
$ws.Range("C135").Value = "La Araucanía"
$ws.Range("D135").Value = 44512
$ws.Range("E135").Value = 9
$ws.Range("F135").Value = 100112040
$ws.Range("G135").Value = "Cilantro"
$ws.Range("H135").Value = "Sin especificar"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 40
$ws.Range("K135").Value = 4000
$ws.Range("L135").Value = 4000
$ws.Range("M135").Value = 4000
$ws.Range("N135").Value = "$/docena de atados (2 kilos)"
$ws.Range("O135").Value = "Provincia de Cautín"
$ws.Range("P135").Value = 2000
$ws.Range("Q135").Value = 2
$ws.Range("R135").Value = "Hortaliza"
